# Insert a new weekly price record at row 15, pushing all existing
# records (old rows 15-164) down by one row (new rows 16-165).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(15).Insert()

$ws.Range("A15").Value = 5
$ws.Range("B15").Value = "Macroferia Regional de Talca"
$ws.Range("C15").Value = "Maule"
$ws.Range("D15").Value = 44685
$ws.Range("E15").Value = 7
$ws.Range("F15").Value = 100112017
$ws.Range("G15").Value = "Apio"
$ws.Range("H15").Value = "Americana (o)"
$ws.Range("I15").Value = "Primera"
$ws.Range("J15").Value = 500
$ws.Range("K15").Value = 8000
$ws.Range("L15").Value = 8000
$ws.Range("M15").Value = 8000
$ws.Range("N15").Value = "$/docena de matas"
$ws.Range("O15").Value = "Provincia del Elquí"
$ws.Range("P15").Value = 1333
$ws.Range("Q15").Value = 6
$ws.Range("R15").Value = "Hortaliza"
